$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume snapshot refresh (GitHub Actions data pull).
# For cells whose new text would otherwise be auto-parsed as a number by
# Excel (losing e.g. trailing zeros), force the cell to Text format just
# for the assignment, then restore the default 'Normal' style so no stray
# formatting is left behind.

$ws.Range('D2').Value = '66.503.72'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '3.183.48'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.64%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.181.71'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('E10').Value = '  -1.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.68'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.515'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.94'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = '3.708.60'
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').Value = '66.545.73'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Value = '3.183.92'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '513.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.56'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.08%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.737'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  -1.08%  '
$ws.Range('E29').Value = '  +7.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.20'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.21'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '514.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0899'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.127'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.76%  '
$ws.Range('D42').Value = '0.0₃0688'
$ws.Range('E42').Value = '  +6.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.301'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.09%  '
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').Value = '2.854.60'
$ws.Range('E46').Value = '  -5.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.54'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('E48').Value = '  +4.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.19%  '
